$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data per the latest scrape.
# Columns B (Coin) and C (Link) are plain text already; set directly.
# Columns D (Price) and E (Volume(1h)) often look numeric to Excel
# (e.g. "238.36", "0.660"), so they are assigned with a leading
# apostrophe to force text, then the style is reset back to Normal
# so no stray text-format style lingers on the cell.

$ws.Range("D2").Value = "'34.872.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.41%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.808.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -3.18%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.16%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'231.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.28%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -1.29%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.17%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'39.34"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -8.26%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +2.30%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0678"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.95%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0992"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.08%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'2.070.83"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -3.13%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.808.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -3.08%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.660"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.48%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'10.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -7.11%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'4.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -4.59%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'34.870.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.47%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'69.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.19%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.0₃0779"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -3.79%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'238.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -4.50%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'11.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -4.75%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -3.10%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.10%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -0.95%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'173.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.05%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'7.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -3.41%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'17.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -4.19%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -3.59%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +6.02%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +0.21%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'3.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.12%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.0546"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.90%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -4.01%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -8.37%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.15"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +5.40%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.683"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'90.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -8.96%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +4.77%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.306.90"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -4.60%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -3.37%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -1.05%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.957"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -6.38%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'14.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -5.61%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'2.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -12.96%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -5.13%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -2.65%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.0510"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.78%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.993.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.09%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.0672"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +7.08%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +0.11%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'98.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -6.16%  "
$ws.Range("E51").Style = "Normal"
